$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C (Weighting), shifting Weighting and the rest right.
$ws.Range("C1").EntireColumn.Insert()
$ws.Range("C1").EntireColumn.ColumnWidth = $ws.Range("B1").EntireColumn.ColumnWidth

# Update header row
$ws.Range("B1").Value = "Reg_method"
$ws.Range("C1").Value = "Cmb_method"

# Fill the new column C with "output_average" for data rows 2-7
$ws.Range("C2:C7").Value = "output_average"

# Local_prec_corr (now column E) values: true/false -> 1/0
$ws.Range("E2").Value = 1
$ws.Range("E3").Value = 0
$ws.Range("E4").Value = 1
$ws.Range("E5").Value = 0
$ws.Range("E6").Value = 1
$ws.Range("E7").Value = 0

# Max_doners (now column F) values: "20" -> 1
$ws.Range("F2:F7").Value = 1

# Run_experiment (now column G) values: true/false -> 1/0
$ws.Range("G2").Value = 1
$ws.Range("G3").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("G5").Value = 0
$ws.Range("G6").Value = 0
$ws.Range("G7").Value = 0

$ws.Range("E8").Select() | Out-Null
